$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.445.22'
$ws.Range('E2').Value = '  -0.46%  '
$ws.Range('D3').Value = '1.568.97'
$ws.Range('E3').Value = '  -1.06%  '
$ws.Range('E4').Value = '  -0.19%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '207.55'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +0.21%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '0.497'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -0.88%  '
$ws.Range('E7').Value = '  -0.16%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '21.97'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -1.16%  '
$ws.Range('E9').Value = '  -1.29%  '
$ws.Range('E10').Value = '  -0.12%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.0866'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  +0.14%  '
$ws.Range('D12').Value = '1.794.86'
$ws.Range('E12').Value = '  -1.02%  '
$ws.Range('D13').Value = '1.564.90'
$ws.Range('E13').Value = '  -2.02%  '
$ws.Range('E14').Value = '  -0.74%  '
$ws.Range('E15').Value = '  -2.84%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '63.28'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  +0.37%  '
$ws.Range('D17').Value = '27.438.91'
$ws.Range('E17').Value = '  -0.50%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '214.02'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -1.36%  '
$ws.Range('E19').Value = '  -0.39%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '7.23'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -0.15%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '4.12'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -0.57%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '9.58'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +0.23%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '2.01'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +1.08%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '153.86'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +0.31%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '6.82'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +1.25%  '
$ws.Range('E27').Value = '  -0.05%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '15.03'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  -0.17%  '
$ws.Range('E29').Value = '  -1.54%  '
$ws.Range('E30').Value = '  -1.06%  '
$ws.Range('E31').Value = '  +1.12%  '
$ws.Range('E32').Value = '  -1.88%  '
$ws.Range('D33').Value = '1.362.49'
$ws.Range('E33').Value = '  -0.30%  '
$ws.Range('E34').Value = '  +0.11%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '1.54'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +1.87%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '0.972'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +0.84%  '
$ws.Range('E37').Value = '  -0.13%  '
$ws.Range('E39').Value = '  -0.91%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.819'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +0.93%  '
$ws.Range('E41').Value = '  -0.09%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.972'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -0.03%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '1.79'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +0.70%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '64.05'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +0.82%  '
$ws.Range('E45').Value = '  -0.08%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '2.15'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -2.03%  '
$ws.Range('D47').Value = '1.704.26'
$ws.Range('E47').Value = '  -1.03%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '85.23'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -2.60%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.0954'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -1.47%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.0495'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -0.43%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +0.03%  '
